# Update generated statistics (F column values) across the four sheets
# as produced by the latest data refresh (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value  = 2644
$ws.Range("F5").Value  = 923
$ws.Range("F6").Value  = 35
$ws.Range("F7").Value  = 2011
$ws.Range("F9").Value  = 203
$ws.Range("F11").Value = 2439
$ws.Range("F12").Value = 528
$ws.Range("F13").Value = 212
$ws.Range("F15").Value = 30
$ws.Range("F16").Value = 115
$ws.Range("F18").Value = 9004
$ws.Range("F20").Value = 6978
$ws.Range("F21").Value = 11402
$ws.Range("F25").Value = 318
$ws.Range("F26").Value = 541
$ws.Range("F27").Value = 2506
$ws.Range("F28").Value = 223
$ws.Range("F29").Value = 189
$ws.Range("F30").Value = 2422
$ws.Range("F31").Value = 610
$ws.Range("F32").Value = 39
$ws.Range("F33").Value = 4483
$ws.Range("F34").Value = 779
$ws.Range("F35").Value = 334
$ws.Range("F37").Value = 495

# --- Sheet "演出" ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F6").Value  = 115
$ws.Range("F8").Value  = 1182

# --- Sheet "本地生活" ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F5").Value  = 137

# --- Sheet "全部类型" ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F6").Value  = 2644
$ws.Range("F8").Value  = 923
$ws.Range("F9").Value  = 35
$ws.Range("F10").Value = 2011
$ws.Range("F14").Value = 203
$ws.Range("F15").Value = 2439
$ws.Range("F16").Value = 115
$ws.Range("F17").Value = 528
$ws.Range("F18").Value = 212
$ws.Range("F21").Value = 115
$ws.Range("F23").Value = 9004
$ws.Range("F25").Value = 6978
$ws.Range("F26").Value = 11402
$ws.Range("F30").Value = 318
$ws.Range("F32").Value = 541
$ws.Range("F34").Value = 2506
$ws.Range("F37").Value = 223
$ws.Range("F38").Value = 189
$ws.Range("F39").Value = 39
$ws.Range("F40").Value = 4483
$ws.Range("F46").Value = 495
